$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 2725.4546
$ws.Range("I86").Value = 2783.3333
$ws.Range("J86").Value = 2656
$ws.Range("K86").Value = 2783.3333
$ws.Range("L86").Value = 2656
$ws.Range("M86").Value = -1660.3333
$ws.Range("N86").Value = -4902
$ws.Range("H89").Value = 2725.4546
$ws.Range("I89").Value = 2783.3333
$ws.Range("J89").Value = 2656
$ws.Range("K89").Value = 13916.6665
$ws.Range("L89").Value = 13280
$ws.Range("M89").Value = -8300.666499999999
$ws.Range("N89").Value = -24512
$ws.Range("H112").Value = 2490.1177
$ws.Range("I112").Value = 625.6
$ws.Range("J112").Value = 3267
$ws.Range("K112").Value = 1876.8
$ws.Range("L112").Value = 9801
$ws.Range("M112").Value = -768.8000000000002
$ws.Range("N112").Value = -12017
$ws.Range("H125").Value = 4309.875
$ws.Range("I125").Value = 3271.8333
$ws.Range("J125").Value = 4932.7
$ws.Range("K125").Value = 29446.4997
$ws.Range("L125").Value = 44394.3
$ws.Range("M125").Value = -26986.4997
$ws.Range("N125").Value = -49314.3
$ws.Range("H127").Value = 1873.55
$ws.Range("I127").Value = 721.5
$ws.Range("J127").Value = 2641.5833
$ws.Range("K127").Value = 2164.5
$ws.Range("L127").Value = 7924.749899999999
$ws.Range("M127").Value = 2795.5
$ws.Range("N127").Value = -17844.7499
$ws.Range("H129").Value = 1173.26
$ws.Range("I129").Value = 366
$ws.Range("J129").Value = 1262.9556
$ws.Range("K129").Value = 1098
$ws.Range("L129").Value = 3788.8668
$ws.Range("M129").Value = 3902
$ws.Range("N129").Value = -13788.8668

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1157.9166
$ws.Range("I74").Value = 705.73334
$ws.Range("J74").Value = 1911.5555
$ws.Range("K74").Value = 705.73334
$ws.Range("L74").Value = 1911.5555
$ws.Range("M74").Value = 168.26666
$ws.Range("N74").Value = -3659.5555
$ws.Range("H77").Value = 1157.9166
$ws.Range("I77").Value = 705.73334
$ws.Range("J77").Value = 1911.5555
$ws.Range("K77").Value = 3528.6667
$ws.Range("L77").Value = 9557.7775
$ws.Range("M77").Value = 839.3333000000002
$ws.Range("N77").Value = -18293.7775
$ws.Range("H113").Value = 37397.332
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 37397.332
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 37397.332
$ws.Range("N113").Value = -46075.332
$ws.Range("H132").Value = 10825.8
$ws.Range("I132").Value = 20726.834
$ws.Range("J132").Value = 4225.1113
$ws.Range("K132").Value = 62180.50199999999
$ws.Range("L132").Value = 12675.3339
$ws.Range("M132").Value = -59650.50199999999
$ws.Range("N132").Value = -17735.3339

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 70000
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 70000
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 70000
$ws.Range("N40").Value = -70530
$ws.Range("H81").Value = 29500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 29500
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 29500
$ws.Range("N81").Value = -31622
$ws.Range("H84").Value = 29500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 29500
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 88500
$ws.Range("N84").Value = -99108
$ws.Range("H86").Value = 256501.5
$ws.Range("I86").Value = 2999.5
$ws.Range("J86").Value = 510003.5
$ws.Range("K86").Value = 2999.5
$ws.Range("L86").Value = 510003.5
$ws.Range("M86").Value = -1876.5
$ws.Range("N86").Value = -512249.5
$ws.Range("H89").Value = 256501.5
$ws.Range("I89").Value = 2999.5
$ws.Range("J89").Value = 510003.5
$ws.Range("K89").Value = 14997.5
$ws.Range("L89").Value = 2550017.5
$ws.Range("M89").Value = -9381.5
$ws.Range("N89").Value = -2561249.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1786.1428
$ws.Range("I33").Value = 575.25
$ws.Range("J33").Value = 3400.6667
$ws.Range("K33").Value = 3451.5
$ws.Range("L33").Value = 20404.0002
$ws.Range("M33").Value = -3168.5
$ws.Range("N33").Value = -20970.0002
$ws.Range("H44").Value = 537.44446
$ws.Range("I44").Value = 309.5
$ws.Range("J44").Value = 719.8
$ws.Range("K44").Value = 928.5
$ws.Range("L44").Value = 2159.4
$ws.Range("M44").Value = -530.5
$ws.Range("N44").Value = -2955.4
$ws.Range("H68").Value = 593
$ws.Range("I68").Value = 640
$ws.Range("J68").Value = 499
$ws.Range("K68").Value = 1920
$ws.Range("L68").Value = 1497
$ws.Range("M68").Value = -1109
$ws.Range("N68").Value = -3119
$ws.Range("H71").Value = 593
$ws.Range("I71").Value = 640
$ws.Range("J71").Value = 499
$ws.Range("K71").Value = 5760
$ws.Range("L71").Value = 4491
$ws.Range("M71").Value = -1704
$ws.Range("N71").Value = -12603
$ws.Range("H94").Value = 3262.8572
$ws.Range("I94").Value = 2000
$ws.Range("J94").Value = 3473.3333
$ws.Range("K94").Value = 6000
$ws.Range("L94").Value = 10419.9999
$ws.Range("M94").Value = -5324
$ws.Range("N94").Value = -11771.9999
$ws.Range("H103").Value = 2757.0527
$ws.Range("I103").Value = 388
$ws.Range("J103").Value = 4139
$ws.Range("K103").Value = 1164
$ws.Range("L103").Value = 12417
$ws.Range("M103").Value = -285
$ws.Range("N103").Value = -14175
$ws.Range("H113").Value = 677.9778
$ws.Range("I113").Value = 598.52
$ws.Range("J113").Value = 777.3
$ws.Range("K113").Value = 1795.56
$ws.Range("L113").Value = 2331.9
$ws.Range("M113").Value = 374.4400000000001
$ws.Range("N113").Value = -6671.9
$ws.Range("H115").Value = 2671.3333
$ws.Range("I115").Value = 2014
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 6042
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -4867
$ws.Range("N115").Value = -11350
$ws.Range("H121").Value = 5061.2856
$ws.Range("I121").Value = 1374
$ws.Range("J121").Value = 7109.778
$ws.Range("K121").Value = 4122
$ws.Range("L121").Value = 21329.334
$ws.Range("M121").Value = -2812
$ws.Range("N121").Value = -23949.334
$ws.Range("H125").Value = 3246.6667
$ws.Range("I125").Value = 530
$ws.Range("J125").Value = 3790
$ws.Range("K125").Value = 1590
$ws.Range("L125").Value = 11370
$ws.Range("M125").Value = 3330
$ws.Range("N125").Value = -21210
$ws.Range("H131").Value = 31396.934
$ws.Range("I131").Value = 372.30768
$ws.Range("J131").Value = 55121.65
$ws.Range("K131").Value = 1116.92304
$ws.Range("L131").Value = 165364.95
$ws.Range("M131").Value = 3923.07696
$ws.Range("N131").Value = -175444.95

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 31000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 31000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 31000
$ws.Range("N62").Value = -32372
$ws.Range("H65").Value = 31000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 31000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 93000
$ws.Range("N65").Value = -99864
$ws.Range("H102").Value = 2588.2666
$ws.Range("I102").Value = 2386.4614
$ws.Range("J102").Value = 3900
$ws.Range("K102").Value = 2386.4614
$ws.Range("L102").Value = 3900
$ws.Range("M102").Value = -764.4614000000001
$ws.Range("N102").Value = -7144

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5700.75
$ws.Range("I7").Value = 3800.5
$ws.Range("J7").Value = 9501.25
$ws.Range("K7").Value = 3800.5
$ws.Range("L7").Value = 9501.25
$ws.Range("M7").Value = -3688.5
$ws.Range("N7").Value = -9725.25
$ws.Range("H61").Value = 102502
$ws.Range("I61").Value = 102502
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 102502
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -102300
$ws.Range("H92").Value = 33000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 33000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 33000
$ws.Range("N92").Value = -37992
$ws.Range("H113").Value = 102502
$ws.Range("I113").Value = 102502
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 102502
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -100332
$ws.Range("H126").Value = 5700.75
$ws.Range("I126").Value = 3800.5
$ws.Range("J126").Value = 9501.25
$ws.Range("K126").Value = 11401.5
$ws.Range("L126").Value = 28503.75
$ws.Range("M126").Value = -8931.5
$ws.Range("N126").Value = -33443.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 168414.33
$ws.Range("I81").Value = 501498.5
$ws.Range("J81").Value = 1872.25
$ws.Range("K81").Value = 1002997
$ws.Range("L81").Value = 3744.5
$ws.Range("M81").Value = -1001936
$ws.Range("N81").Value = -5866.5
$ws.Range("H84").Value = 168414.33
$ws.Range("I84").Value = 501498.5
$ws.Range("J84").Value = 1872.25
$ws.Range("K84").Value = 5014985
$ws.Range("L84").Value = 18722.5
$ws.Range("M84").Value = -5009681
$ws.Range("N84").Value = -29330.5
$ws.Range("H122").Value = 3035.037
$ws.Range("I122").Value = 2234.8823
$ws.Range("J122").Value = 4395.3
$ws.Range("K122").Value = 6704.646900000001
$ws.Range("L122").Value = 13185.9
$ws.Range("M122").Value = -4254.646900000001
$ws.Range("N122").Value = -18085.9
$ws.Range("H126").Value = 9859.632
$ws.Range("I126").Value = 14030.818
$ws.Range("J126").Value = 4124.25
$ws.Range("K126").Value = 42092.454
$ws.Range("L126").Value = 12372.75
$ws.Range("M126").Value = -39622.454
$ws.Range("N126").Value = -17312.75
